$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames (row 1)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B7').Value = 'Rincón De Romos'
$ws.Range('B25').Value = 'Comitán De Domínguez'
$ws.Range('B34').Value = 'Mazapa De Madero'
$ws.Range('B55').Value = 'Hidalgo Del Parral'
$ws.Range('B66').Value = 'San Francisco De Borja'
$ws.Range('B86').Value = 'Villa De Álvarez'
$ws.Range('A88').Value = 'Ciudad De México'
$ws.Range('B115').Value = 'San Juan De Guadalupe'
$ws.Range('A120').Value = 'Estado De México'
$ws.Range('B120').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B121').Value = 'Almoloya De Alquisiras'
$ws.Range('B122').Value = 'Almoloya De Juárez'
$ws.Range('B125').Value = 'Atizapán De Zaragoza'
$ws.Range('B129').Value = 'Coacalco De Berriozábal'
$ws.Range('B132').Value = 'Ecatepec De Morelos'
$ws.Range('B134').Value = 'Ixtapan De La Sal'
$ws.Range('B140').Value = 'Naucalpan De Juárez'
$ws.Range('B146').Value = 'San Martín De Las Pirámides'
$ws.Range('B147').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B153').Value = 'Tenango Del Valle'
$ws.Range('B155').Value = 'Tlalnepantla De Baz'
$ws.Range('B165').Value = 'San Miguel De Allende'
$ws.Range('B166').Value = 'Apaseo El Alto'
$ws.Range('B171').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B174').Value = 'Jaral Del Progreso'
$ws.Range('B181').Value = 'Purísima Del Rincón'
$ws.Range('B186').Value = 'San Francisco Del Rincón'
$ws.Range('B188').Value = 'San Luis De La Paz'
$ws.Range('B190').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B192').Value = 'Silao De La Victoria'
$ws.Range('B196').Value = 'Valle De Santiago'
$ws.Range('B201').Value = 'Acapulco De Juárez'
$ws.Range('B202').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B204').Value = 'Atenango Del Río'
$ws.Range('B207').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B210').Value = 'Coyuca De Catalán'
$ws.Range('B213').Value = 'Cuetzala Del Progreso'
$ws.Range('B216').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B217').Value = 'Iguala De La Independencia'
$ws.Range('B218').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B219').Value = 'Zihuatanejo De Azueta'
$ws.Range('B228').Value = 'Taxco De Alarcón'
$ws.Range('B229').Value = 'Técpan De Galeana'
$ws.Range('B231').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B238').Value = 'Cuautepec De Hinojosa'
$ws.Range('B239').Value = 'Huejutla De Reyes'
$ws.Range('B242').Value = 'Jacala De Ledezma'
$ws.Range('B245').Value = 'Mineral Del Chico'
$ws.Range('B246').Value = 'Molango De Escamilla'
$ws.Range('B247').Value = 'Pachuca De Soto'
$ws.Range('B249').Value = 'Tenango De Doria'
$ws.Range('B250').Value = 'Tepehuacán De Guerrero'
$ws.Range('B251').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B254').Value = 'Tulancingo De Bravo'
$ws.Range('B255').Value = 'Zacualtipán De Ángeles'
$ws.Range('B259').Value = 'Ahualulco De Mercado'
$ws.Range('B263').Value = 'Atemajac De Brizuela'
$ws.Range('B264').Value = 'Atotonilco El Alto'
$ws.Range('B265').Value = 'Autlán De Navarro'
$ws.Range('B275').Value = 'Encarnación De Díaz'
$ws.Range('B279').Value = 'Huejuquilla El Alto'
$ws.Range('B285').Value = 'Lagos De Moreno'
$ws.Range('B288').Value = 'Ojuelos De Jalisco'
$ws.Range('B293').Value = 'San Juan De Los Lagos'
$ws.Range('B294').Value = 'San Martín De Bolaños'
$ws.Range('B296').Value = 'San Miguel El Alto'
$ws.Range('B297').Value = 'San Sebastián Del Oeste'
$ws.Range('B299').Value = 'Tamazula De Gordiano'
$ws.Range('B301').Value = 'Teocuitatlán De Corona'
$ws.Range('B302').Value = 'Tepatitlán De Morelos'
$ws.Range('B303').Value = 'Tizapán El Alto'
$ws.Range('B310').Value = 'Unión De Tula'
$ws.Range('B314').Value = 'Yahualica De González Gallo'
$ws.Range('B315').Value = 'Zacoalco De Torres'
$ws.Range('B318').Value = 'Zapotitlán De Vadillo'
$ws.Range('B319').Value = 'Zapotlán Del Rey'
$ws.Range('B320').Value = 'Zapotlán El Grande'
$ws.Range('B338').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B340').Value = 'Cojumatlán De Régules'
$ws.Range('B394').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B416').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B421').Value = 'Tetela Del Volcán'
$ws.Range('B429').Value = 'Amatlán De Cañas'
$ws.Range('B433').Value = 'Ixtlán Del Río'
$ws.Range('B439').Value = 'Santa María Del Oro'
$ws.Range('B450').Value = 'Lampazos De Naranjo'
$ws.Range('B453').Value = 'San Nicolás De Los Garza'
$ws.Range('B456').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B460').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B461').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B462').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B465').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B466').Value = 'Nejapa De Madero'
$ws.Range('B467').Value = 'Oaxaca De Juárez'
$ws.Range('B468').Value = 'Ocotlán De Morelos'
$ws.Range('B470').Value = 'Putla Villa De Guerrero'
$ws.Range('B485').Value = 'San Mateo Del Mar'
$ws.Range('B505').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B506').Value = 'Teotitlán De Flores Magón'
$ws.Range('B507').Value = 'Tlacolula De Matamoros'
$ws.Range('B508').Value = 'Tlalixtac De Cabrera'
$ws.Range('B509').Value = 'Villa De Tututepec'
$ws.Range('B510').Value = 'Villa Sola De Vega'
$ws.Range('B511').Value = 'Villa Talea De Castro'
$ws.Range('B512').Value = 'Zimatlán De Álvarez'
$ws.Range('B518').Value = 'Ayotoxco De Guerrero'
$ws.Range('B536').Value = 'Tepanco De López'
$ws.Range('B539').Value = 'Tepexi De Rodríguez'
$ws.Range('B540').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B546').Value = 'Amealco De Bonfil'
$ws.Range('B547').Value = 'Cadereyta De Montes'
$ws.Range('B549').Value = 'Jalpan De Serra'
$ws.Range('B552').Value = 'Pinal De Amoles'
$ws.Range('B555').Value = 'San Juan Del Río'
$ws.Range('B558').Value = 'Armadillo De Los Infante'
$ws.Range('B560').Value = 'Ciudad Del Maíz'
$ws.Range('B563').Value = 'Mexquitic De Carmona'
$ws.Range('B570').Value = 'Villa De Ramos'
$ws.Range('B626').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B635').Value = 'Amatlán De Los Reyes'
$ws.Range('B639').Value = 'Cosamaloapan De Carpio'
$ws.Range('B641').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B646').Value = 'Lerdo De Tejada'
$ws.Range('B653').Value = 'Poza Rica De Hidalgo'
$ws.Range('B676').Value = 'Concepción Del Oro'
$ws.Range('B678').Value = 'El Plateado De Joaquín Amaro'
$ws.Range('B689').Value = 'Nochistlán De Mejía'
$ws.Range('B697').Value = 'Teúl De González Ortega'
$ws.Range('B698').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B700').Value = 'Villa De Cos'

# Remove trailing metadata rows 708-712 in one shot (shifts rows up, trims dimension to D706)
$ws.Range("A708:A712").EntireRow.Delete()
